# Scheduled market-data refresh: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for affected leve rows across the crafting-profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 688.4
$ws.Range("I33").Value = 757.375
$ws.Range("J33").Value = 412.5
$ws.Range("K33").Value = 757.375
$ws.Range("L33").Value = 412.5
$ws.Range("M33").Value = -528.375
$ws.Range("N33").Value = -870.5

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 3156.1052
$ws.Range("I40").Value = 2244.5
$ws.Range("J40").Value = 3399.2
$ws.Range("K40").Value = 2244.5
$ws.Range("L40").Value = 3399.2
$ws.Range("M40").Value = -2069.5
$ws.Range("N40").Value = -3749.2

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 2976
$ws.Range("I100").Value = 2290
$ws.Range("K100").Value = 2290
$ws.Range("M100").Value = -1749

# Row 129: Practical Command
$ws.Range("H129").Value = 1418.0416
$ws.Range("J129").Value = 1762.5555
$ws.Range("L129").Value = 5287.666499999999
$ws.Range("N129").Value = -15287.6665

$ws = $wb.Worksheets.Item("ARM")
# Row 58: Some Dragoons Have All the Luck
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20860

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 3085.4348
$ws.Range("I63").Value = 2089.7368
$ws.Range("J63").Value = 7815
$ws.Range("K63").Value = 2089.7368
$ws.Range("L63").Value = 7815
$ws.Range("M63").Value = -1403.7368
$ws.Range("N63").Value = -9187

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 3085.4348
$ws.Range("I66").Value = 2089.7368
$ws.Range("J66").Value = 7815
$ws.Range("K66").Value = 10448.684
$ws.Range("L66").Value = 39075
$ws.Range("M66").Value = -7016.684000000001
$ws.Range("N66").Value = -45939

# Row 97: Ore for Me
$ws.Range("H97").Value = 1125
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -504
$ws.Range("N97").Value = -2242

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1700.8
$ws.Range("I102").Value = 1634.2222
$ws.Range("K102").Value = 1634.2222
$ws.Range("M102").Value = -12.22219999999993

# Row 131: Additions to the Armoire
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1516.6666
$ws.Range("I86").Value = 1510
$ws.Range("J86").Value = 1525
$ws.Range("K86").Value = 1510
$ws.Range("L86").Value = 1525
$ws.Range("M86").Value = -387
$ws.Range("N86").Value = -3771

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1516.6666
$ws.Range("I89").Value = 1510
$ws.Range("J89").Value = 1525
$ws.Range("K89").Value = 7550
$ws.Range("L89").Value = 7625
$ws.Range("M89").Value = -1934
$ws.Range("N89").Value = -18857

# Row 94: High Steal
$ws.Range("H94").Value = 1346.091
$ws.Range("I94").Value = 1203
$ws.Range("J94").Value = 1990
$ws.Range("K94").Value = 1203
$ws.Range("L94").Value = 1990
$ws.Range("M94").Value = -752
$ws.Range("N94").Value = -2892

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2462.5
$ws.Range("I99").Value = 1535.7142
$ws.Range("K99").Value = 1535.7142
$ws.Range("M99").Value = -37.71419999999989

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3826
$ws.Range("I105").Value = 3782.5
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3782.5
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -2035.5
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 6291444.5
$ws.Range("I31").Value = 1628.1111
$ws.Range("J31").Value = 19611056
$ws.Range("K31").Value = 1628.1111
$ws.Range("L31").Value = 19611056
$ws.Range("M31").Value = -1333.1111
$ws.Range("N31").Value = -19611646

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 6291444.5
$ws.Range("I34").Value = 1628.1111
$ws.Range("J34").Value = 19611056
$ws.Range("K34").Value = 1628.1111
$ws.Range("L34").Value = 19611056
$ws.Range("M34").Value = -1426.1111
$ws.Range("N34").Value = -19611460

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 760432.9
$ws.Range("I58").Value = 1358.3334
$ws.Range("J58").Value = 2088813.2
$ws.Range("K58").Value = 1358.3334
$ws.Range("L58").Value = 2088813.2
$ws.Range("M58").Value = -1155.3334
$ws.Range("N58").Value = -2089219.2

# Row 136: Turali Quality
$ws.Range("H136").Value = 760432.9
$ws.Range("I136").Value = 1358.3334
$ws.Range("J136").Value = 2088813.2
$ws.Range("K136").Value = 4075.0002
$ws.Range("L136").Value = 6266439.6
$ws.Range("M136").Value = -1525.0002
$ws.Range("N136").Value = -6271539.6

$ws = $wb.Worksheets.Item("CUL")
# Row 13: Fishy Revelations
$ws.Range("H13").Value = 2072
$ws.Range("I13").Value = 900.3333
$ws.Range("J13").Value = 2574.1428
$ws.Range("K13").Value = 2700.9999
$ws.Range("L13").Value = 7722.428400000001
$ws.Range("M13").Value = -2532.9999
$ws.Range("N13").Value = -8058.428400000001

# Row 58: Bread in the Clouds
$ws.Range("H58").Value = 2675
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15256

# Row 132: More Mezcal
$ws.Range("H132").Value = 3444.889
$ws.Range("I132").Value = 1066.6666
$ws.Range("J132").Value = 4634
$ws.Range("K132").Value = 9599.999400000001
$ws.Range("L132").Value = 41706
$ws.Range("M132").Value = -7069.999400000001
$ws.Range("N132").Value = -46766

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 1923.8235
$ws.Range("I80").Value = 2040.4
$ws.Range("J80").Value = 1600
$ws.Range("K80").Value = 2040.4
$ws.Range("L80").Value = 1600
$ws.Range("M80").Value = -1042.4
$ws.Range("N80").Value = -3596

# Row 82: Appeasing the Astromancer
$ws.Range("H82").Value = 45560
$ws.Range("J82").Value = 45560
$ws.Range("L82").Value = 45560
$ws.Range("N82").Value = -46326

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1923.8235
$ws.Range("I83").Value = 2040.4
$ws.Range("J83").Value = 1600
$ws.Range("K83").Value = 10202
$ws.Range("L83").Value = 8000
$ws.Range("M83").Value = -5210
$ws.Range("N83").Value = -17984

# Row 85: Silver Bar of Upcycling (L)
$ws.Range("H85").Value = 45560
$ws.Range("J85").Value = 45560
$ws.Range("L85").Value = 45560
$ws.Range("N85").Value = -48212

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 2911.2
$ws.Range("I97").Value = 2670.25
$ws.Range("K97").Value = 2670.25
$ws.Range("M97").Value = -2174.25

# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 626.16
$ws.Range("I107").Value = 585.0625
$ws.Range("J107").Value = 699.2222
$ws.Range("K107").Value = 585.0625
$ws.Range("L107").Value = 699.2222
$ws.Range("M107").Value = 1334.9375
$ws.Range("N107").Value = -4539.2222

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3333880
$ws.Range("I46").Value = 825
$ws.Range("J46").Value = 9999990
$ws.Range("K46").Value = 825
$ws.Range("L46").Value = 9999990
$ws.Range("M46").Value = -637
$ws.Range("N46").Value = -10000366

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 424.29413
$ws.Range("I55").Value = 270.07693
$ws.Range("J55").Value = 925.5
$ws.Range("K55").Value = 270.07693
$ws.Range("L55").Value = 925.5
$ws.Range("M55").Value = -97.07693
$ws.Range("N55").Value = -1271.5

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 1302.125
$ws.Range("I82").Value = 572.9231
$ws.Range("J82").Value = 2163.9092
$ws.Range("K82").Value = 572.9231
$ws.Range("L82").Value = 2163.9092
$ws.Range("M82").Value = -211.9231
$ws.Range("N82").Value = -2885.9092

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 1302.125
$ws.Range("I85").Value = 572.9231
$ws.Range("J85").Value = 2163.9092
$ws.Range("K85").Value = 572.9231
$ws.Range("L85").Value = 2163.9092
$ws.Range("M85").Value = 675.0769
$ws.Range("N85").Value = -4659.9092

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 2397.08
$ws.Range("I93").Value = 2171.7368
$ws.Range("J93").Value = 3110.6667
$ws.Range("K93").Value = 2171.7368
$ws.Range("L93").Value = 3110.6667
$ws.Range("M93").Value = -923.7368000000001
$ws.Range("N93").Value = -5606.6667

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 3013.6428
$ws.Range("I100").Value = 2562.875
$ws.Range("K100").Value = 2562.875
$ws.Range("M100").Value = -2021.875
